# Update gh-pages output data (regenerated counts) on sheets "展览" (1) and "全部类型" (4).
# Both sheets contain identical tables, so the same set of edits is applied to each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value  = 1314
    $ws.Range("F4").Value  = 144
    $ws.Range("F10").Value = 126
    $ws.Range("F11").Value = 4428
    $ws.Range("F12").Value = 6697
    $ws.Range("F14").Value = 54
    $ws.Range("G15").Value = "不可售"
    $ws.Range("F16").Value = 563
    $ws.Range("F18").Value = 4097
    $ws.Range("F19").Value = 460
    $ws.Range("F20").Value = 67
    $ws.Range("F22").Value = 2675
    $ws.Range("F24").Value = 544
    $ws.Range("F26").Value = 342
    $ws.Range("F27").Value = 346
    $ws.Range("F28").Value = 393
    $ws.Range("F31").Value = 1610
    $ws.Range("F32").Value = 1013
    $ws.Range("F33").Value = 58
    $ws.Range("F34").Value = 120
    $ws.Range("F36").Value = 534
    $ws.Range("F40").Value = 622
}

$wb.Save()
